$wb = $excel.ActiveWorkbook

# --- "Data" sheet: append the newest weekly observation ---
$wsData = $wb.Worksheets.Item("Data")

# Copy formatting (number format / font / border / alignment) from the last
# existing data row down into the new row, then fill in the new values.
$wsData.Range("A94:B94").Copy($wsData.Range("A95:B95"))
$wsData.Range("A95").Value = 45126
$wsData.Range("B95").Value = 3230.457

# --- "SeriesInfo" sheet: refresh the metadata pulled from the FRED API ---
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end
$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Range("B3").Value = "2023-07-24"
$wsInfo.Range("B3").Style = "Normal"

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Range("B4").Value = "2023-07-24"
$wsInfo.Range("B4").Style = "Normal"

# observation_end
$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Range("B7").Value = "2023-07-19"
$wsInfo.Range("B7").Style = "Normal"

# last_updated
$wsInfo.Range("B14").NumberFormat = "@"
$wsInfo.Range("B14").Value = "2023-07-20 15:35:27-05"
$wsInfo.Range("B14").Style = "Normal"

# popularity
$wsInfo.Range("B15").Value = 78
